# Edit jogos.xlsx:
#  1. Rename the "PC" platform entries to "PC - Steam" (most rows) or
#     "PC - Xbox " (rows that correspond to games bought on the Xbox/PC
#     Game Pass storefront: "Lies of P" and "Hollow Knight" PC entries).
#  2. Fix the "Naruto Shippuden" title by removing the colon.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Plataforma" (column C) should become "PC - Xbox " instead of
# "PC - Steam".
$xboxRows = @(73, 74)

$lastRow = 76
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Text -eq "PC") {
        if ($xboxRows -contains $r) {
            $cell.Value = "PC - Xbox "
        } else {
            $cell.Value = "PC - Steam"
        }
    }
}

# Fix game title spelling (remove colon after "Shippuden").
$ws.Cells.Item(40, 1).Value = "Naruto Shippuden Ultimate Ninja Storm 4"

# Reflect where the user was last working in the sheet (cosmetic view state).
$win = $excel.ActiveWindow
$ws.Range("C74").Select() | Out-Null
$win.ScrollRow = 31
$win.ScrollColumn = 1
